# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# "175c5e52..." file has now been handed back (in sync with en-US), while the
# "8be37c0f..." file remains in translation. Row 2 / Row 3 swap identity on
# every sheet, and the per-language sheets gain a populated "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" set of columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "'175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$overview.Range("C2").Value = ".md"
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("G2").Value = "'2016-09-03 08:31:33"

$overview.Range("A3").Value = "'8be37c0f-042b-48dd-b810-b5e49c366266ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$overview.Range("C3").Value = ".md"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("G3").Value = "'2016-09-03 08:30:48"

$overview.Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e80a0ace2511f16f40f71139be338c78456157f1/e2e/8be37c0f-042b-48dd-b810-b5e49c366266ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md", [System.Type]::Missing, [System.Type]::Missing, "'e2e\175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")
$overview.Hyperlinks.Add($overview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3741608aec8ce568748b07ad085e6e23c75d2a57/e2e/175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md", [System.Type]::Missing, [System.Type]::Missing, "'e2e\8be37c0f-042b-48dd-b810-b5e49c366266ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")

$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "'175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$zhcn.Range("B2").Value = ".md"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("D2").Value = "e2e"
$zhcn.Range("E2").Value = "ht"
$zhcn.Range("F2").Value = "'False"
$zhcn.Range("G2").Value = "'175c5e52-4b3b-4df7-b275-7406c7a6b704oooooooooooooooooooooooooooooooooooooooo.0e2600ce9f4c4300e6d0a0806bca16914d3f1300.zh-cn.xlf"
$zhcn.Range("H2").Value = "'2016-09-03 08:31:29"
$zhcn.Range("I2").Value = "'175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$zhcn.Range("J2").Value = "'175c5e52-4b3b-4df7-b275-7406c7a6b704oooooooooooooooooooooooooooooooooooooooo.0e2600ce9f4c4300e6d0a0806bca16914d3f1300.zh-cn.xlf"
$zhcn.Range("K2").Value = "'2016-09-03 08:31:45"
$zhcn.Range("L2").Value = "'"
$zhcn.Range("M2").Value = "'True"
$zhcn.Range("N2").Value = "'"
$zhcn.Range("O2").Value = "'False"
$zhcn.Range("P2").Value = "'"

$zhcn.Range("A3").Value = "'8be37c0f-042b-48dd-b810-b5e49c366266ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$zhcn.Range("B3").Value = ".md"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("D3").Value = "e2e"
$zhcn.Range("E3").Value = "ht"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "'8be37c0f-042b-48dd-b810-b5e49c366266oooooooooooooooooooooooooooooooooooooooo.5f8ba311eede70b4fd2894c9c520ea9b9d57a630.zh-cn.xlf"
$zhcn.Range("H3").Value = "'2016-09-03 08:30:43"
$zhcn.Range("I3").Value = "'"
$zhcn.Range("J3").Value = "'"
$zhcn.Range("K3").Value = "'0001-01-01 00:00:00"
$zhcn.Range("L3").Value = "'"
$zhcn.Range("M3").Value = "'True"
$zhcn.Range("N3").Value = "'"
$zhcn.Range("O3").Value = "'False"
$zhcn.Range("P3").Value = "'"

$zhcn.Range("I2").Style = "Hyperlink"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3741608aec8ce568748b07ad085e6e23c75d2a57/e2e/175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md", [System.Type]::Missing, [System.Type]::Missing, "'175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3741608aec8ce568748b07ad085e6e23c75d2a57/e2e/175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md", [System.Type]::Missing, [System.Type]::Missing, "'175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e80a0ace2511f16f40f71139be338c78456157f1/e2e/8be37c0f-042b-48dd-b810-b5e49c366266ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md", [System.Type]::Missing, [System.Type]::Missing, "'8be37c0f-042b-48dd-b810-b5e49c366266ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "'175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$dede.Range("B2").Value = ".md"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("D2").Value = "e2e"
$dede.Range("E2").Value = "ht"
$dede.Range("F2").Value = "'False"
$dede.Range("G2").Value = "'175c5e52-4b3b-4df7-b275-7406c7a6b704oooooooooooooooooooooooooooooooooooooooo.0e2600ce9f4c4300e6d0a0806bca16914d3f1300.de-de.xlf"
$dede.Range("H2").Value = "'2016-09-03 08:31:33"
$dede.Range("I2").Value = "'175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$dede.Range("J2").Value = "'175c5e52-4b3b-4df7-b275-7406c7a6b704oooooooooooooooooooooooooooooooooooooooo.0e2600ce9f4c4300e6d0a0806bca16914d3f1300.de-de.xlf"
$dede.Range("K2").Value = "'2016-09-03 08:31:51"
$dede.Range("L2").Value = "'"
$dede.Range("M2").Value = "'True"
$dede.Range("N2").Value = "'"
$dede.Range("O2").Value = "'False"
$dede.Range("P2").Value = "'"

$dede.Range("A3").Value = "'8be37c0f-042b-48dd-b810-b5e49c366266ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = "In Translation"
$dede.Range("D3").Value = "e2e"
$dede.Range("E3").Value = "ht"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "'8be37c0f-042b-48dd-b810-b5e49c366266oooooooooooooooooooooooooooooooooooooooo.5f8ba311eede70b4fd2894c9c520ea9b9d57a630.de-de.xlf"
$dede.Range("H3").Value = "'2016-09-03 08:30:48"
$dede.Range("I3").Value = "'"
$dede.Range("J3").Value = "'"
$dede.Range("K3").Value = "'0001-01-01 00:00:00"
$dede.Range("L3").Value = "'"
$dede.Range("M3").Value = "'True"
$dede.Range("N3").Value = "'"
$dede.Range("O3").Value = "'False"
$dede.Range("P3").Value = "'"

$dede.Range("I2").Style = "Hyperlink"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3741608aec8ce568748b07ad085e6e23c75d2a57/e2e/175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md", [System.Type]::Missing, [System.Type]::Missing, "'175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3741608aec8ce568748b07ad085e6e23c75d2a57/e2e/175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md", [System.Type]::Missing, [System.Type]::Missing, "'175c5e52-4b3b-4df7-b275-7406c7a6b704ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e80a0ace2511f16f40f71139be338c78456157f1/e2e/8be37c0f-042b-48dd-b810-b5e49c366266ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md", [System.Type]::Missing, [System.Type]::Missing, "'8be37c0f-042b-48dd-b810-b5e49c366266ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40
